$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 446, shifting all rows 446:473 down to 447:474
$ws.Rows(446).Insert()

# Populate the newly inserted row 446 with the latest week's price report
$ws.Range("A446").Value = 7
$ws.Range("B446").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C446").Value = "Ñuble"
$ws.Range("D446").Value = 44931
$ws.Range("E446").Value = 16
$ws.Range("F446").Value = 100114001
$ws.Range("G446").Value = "Papa"
$ws.Range("H446").Value = "Patagonia"
$ws.Range("I446").Value = "1a nueva(o)"
$ws.Range("J446").Value = 200
$ws.Range("K446").Value = 12000
$ws.Range("L446").Value = 13000
$ws.Range("M446").Value = 12500
$ws.Range("N446").Value = "$/saco 25 kilos"
$ws.Range("O446").Value = "Región de La Araucanía"
$ws.Range("P446").Value = 500
$ws.Range("Q446").Value = 25
$ws.Range("R446").Value = "Hortaliza"
